# Refresh the Price (D) and Volume(1h) (E) columns of the cryptos
# list with the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text (e.g. "67.407.82", "8.00").
# A leading apostrophe forces Excel to keep the entry as text
# instead of reinterpreting it as a number (which would collapse
# "8.00" to 8 and drop the cell's General number format).
$ws.Range("D2").Value = "'67.407.82"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "'3.527.20"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'595.89"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'173.92"
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.84%  "
$ws.Range("E9").Value = "  +7.12%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "'4.137.36"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("E14").Value = "  +1.79%  "
$ws.Range("D15").Value = "'0.0000183"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "'67.315.63"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "'3.527.33"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").Value = "'398.12"
$ws.Range("E20").Value = "  +1.65%  "
$ws.Range("D21").Value = "'8.00"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").Value = "'73.55"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'0.0000123"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("D26").Value = "'10.32"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Value = "'0.180"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").Value = "'24.13"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'1.64"
$ws.Range("E34").Value = "  +3.61%  "
$ws.Range("D35").Value = "'164.08"
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("D36").Value = "'0.897"
$ws.Range("E36").Value = "  -1.75%  "
$ws.Range("E37").Value = "  -1.06%  "
$ws.Range("D38").Value = "'6.96"
$ws.Range("E38").Value = "  +3.64%  "
$ws.Range("D39").Value = "'4.73"
$ws.Range("E39").Value = "  +1.53%  "
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("D42").Value = "'26.56"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("D44").Value = "'2.804.05"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").Value = "'42.98"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").Value = "'0.0311"
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").Value = "'342.42"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("D48").Value = "'1.11"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("D49").Value = "'33.93"
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "'6.55"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("E51").Value = "  +0.11%  "
